$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.798.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.889.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4785"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2962"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06634"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "100.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +18.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.883.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07552"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.150"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6618"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "302.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +25.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.786.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.43%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007595"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.133.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.157"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.221"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.327"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.956"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1128"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.353"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.176"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05087"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.165"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7337"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.718"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("E37").Value = "  +3.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.706"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.062"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8972"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "108.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4202"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.656"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.374"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.135"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1231"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05645"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.400"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.53%  "
